$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 45845
$ws.Range("A3").Value = 251310
$ws.Range("B3").Value = 45770
$ws.Range("D3").Value = 45769
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 8611
$ws.Range("G3").Value = 336
$ws.Range("H3").Value = "'6"
$ws.Range("H3").ClearFormats()
$ws.Range("J3").Value = 410
$ws.Range("K3").Value = 820
$ws.Range("A4").Value = 251905
$ws.Range("D2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "CAMPO VUOTO"
$ws.Range("C4").Value = "IN STAMPA"
$ws.Range("D4").Value = 45805
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 13841
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = "CAMPO VUOTO"
$ws.Range("J4").Value = 165
$ws.Range("K4").Value = 825
$ws.Range("A1").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 76
$ws.Range("N4").Value = 0
$ws.Range("A5").Value = 252666
$ws.Range("D2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = "CAMPO VUOTO"
$ws.Range("C5").Value = "DA STAMPARE"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 13175
$ws.Range("G5").Value = 600
$ws.Range("H5").Value = "CAMPO VUOTO"
$ws.Range("J5").Value = 680
$ws.Range("K5").Value = 680
$ws.Range("A1").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 70
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = "CAMPO VUOTO"
$ws.Range("P5").Value = "CAMPO VUOTO"
$ws.Range("A6").Value = 251704
$ws.Range("B6").Value = 45846
$ws.Range("B2").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = 45848
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 28904
$ws.Range("G6").Value = 2217
$ws.Range("H6").Value = "'1"
$ws.Range("H6").ClearFormats()
$ws.Range("I6").Value = "bobina"
$ws.Range("J6").Value = 340
$ws.Range("K6").Value = 1020
$ws.Range("D2").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = "CAMPO VUOTO"
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = "X"
$ws.Range("P6").Value = 40054
$ws.Range("A7").Value = 252453
$ws.Range("B7").Value = 45845
$ws.Range("D7").Value = 45835
$ws.Range("F7").Value = 2990
$ws.Range("G7").Value = 169
$ws.Range("H7").Value = "3 / 4"
$ws.Range("D2").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "CAMPO VUOTO"
$ws.Range("J7").Value = 190
$ws.Range("K7").Value = 760
$ws.Range("M7").Value = " ERR tipologia incarto non definita"
$ws.Range("O7").Value = "CAMPO VUOTO"
$ws.Range("P7").Value = "CAMPO VUOTO"
$ws.Range("A8").Value = 252983
$ws.Range("B8").Value = 45847
$ws.Range("D8").Value = 45862
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1912
$ws.Range("G8").Value = 97
$ws.Range("H8").Value = "'3"
$ws.Range("H8").ClearFormats()
$ws.Range("J8").Value = 85
$ws.Range("K8").Value = 425
$ws.Range("N8").Value = 4
$ws.Range("P8").Value = 40055
$ws.Range("A9").Value = 252450
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 45852
$ws.Range("C9").Value = "IN STAMPA"
$ws.Range("B2").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = 45849
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 5718
$ws.Range("G9").Value = 300
$ws.Range("H9").Value = "'6"
$ws.Range("H9").ClearFormats()
$ws.Range("I9").Value = "bobina"
$ws.Range("J9").Value = 390
$ws.Range("K9").Value = 780
$ws.Range("D2").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = "CAMPO VUOTO"
$ws.Range("M9").Value = "Dati OK"
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = "CAMPO VUOTO"
$ws.Range("P9").Value = "CAMPO VUOTO"
$ws.Range("A10").Value = 252354
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = 45849
$ws.Range("C10").Value = "IN STAMPA"
$ws.Range("B2").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Value = 45847
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 5060
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = "'1"
$ws.Range("H10").ClearFormats()
$ws.Range("I10").Value = "bobina"
$ws.Range("J10").Value = 85
$ws.Range("K10").Value = 595
$ws.Range("D2").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Range("L10").Value = "CAMPO VUOTO"
$ws.Range("M10").Value = "Dati OK"
$ws.Range("N10").Value = 2
$ws.Range("O10").Value = "CAMPO VUOTO"
$ws.Range("P10").Value = "CAMPO VUOTO"

Write-Host "Update complete"
